# Auto-generated script to apply numeric value updates to the Tonberry Profits workbook
# Sheets are FFXIV Disciple of the Hand crafting-leve profit trackers: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4760.5
$ws.Range("I76").Value = 5146.25
$ws.Range("J76").Value = 4374.75
$ws.Range("K76").Value = 5146.25
$ws.Range("L76").Value = 4374.75
$ws.Range("M76").Value = -4831.25
$ws.Range("N76").Value = -5004.75
$ws.Range("H79").Value = 4760.5
$ws.Range("I79").Value = 5146.25
$ws.Range("J79").Value = 4374.75
$ws.Range("K79").Value = 5146.25
$ws.Range("L79").Value = 4374.75
$ws.Range("M79").Value = -4054.25
$ws.Range("N79").Value = -6558.75
$ws.Range("H80").Value = 2643.4285
$ws.Range("I80").Value = 4401
$ws.Range("J80").Value = 300
$ws.Range("K80").Value = 13203
$ws.Range("L80").Value = 900
$ws.Range("M80").Value = -12205
$ws.Range("N80").Value = -2896
$ws.Range("H83").Value = 2643.4285
$ws.Range("I83").Value = 4401
$ws.Range("J83").Value = 300
$ws.Range("K83").Value = 39609
$ws.Range("L83").Value = 2700
$ws.Range("M83").Value = -34617
$ws.Range("N83").Value = -12684
$ws.Range("H99").Value = 1819280.6
$ws.Range("I99").Value = 312.83334
$ws.Range("K99").Value = 938.5000200000001
$ws.Range("M99").Value = 559.4999799999999
$ws.Range("H132").Value = 1214.3914
$ws.Range("I132").Value = 1169.591
$ws.Range("K132").Value = 3508.773
$ws.Range("M132").Value = -978.7729999999997
$ws.Range("H138").Value = 3399.2354
$ws.Range("J138").Value = 2954.3438
$ws.Range("L138").Value = 8863.0314
$ws.Range("N138").Value = -19143.0314

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13315.5
$ws.Range("I32").Value = 8375.484
$ws.Range("J32").Value = 19836.32
$ws.Range("K32").Value = 8375.484
$ws.Range("L32").Value = 19836.32
$ws.Range("M32").Value = -8088.484
$ws.Range("N32").Value = -20410.32
$ws.Range("H44").Value = 22750
$ws.Range("H55").Value = 18093.334
$ws.Range("J55").Value = 18093.334
$ws.Range("L55").Value = 18093.334
$ws.Range("N55").Value = -18723.334
$ws.Range("H74").Value = 815.0952
$ws.Range("I74").Value = 590.2162
$ws.Range("K74").Value = 590.2162
$ws.Range("M74").Value = 283.7838
$ws.Range("H77").Value = 815.0952
$ws.Range("I77").Value = 590.2162
$ws.Range("K77").Value = 2951.081
$ws.Range("M77").Value = 1416.919
$ws.Range("H80").Value = 49500
$ws.Range("J80").Value = 49500
$ws.Range("L80").Value = 49500
$ws.Range("N80").Value = -51496
$ws.Range("H83").Value = 49500
$ws.Range("J83").Value = 49500
$ws.Range("L83").Value = 148500
$ws.Range("N83").Value = -158484

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 49800
$ws.Range("J57").Value = 49800
$ws.Range("L57").Value = 49800
$ws.Range("N57").Value = -51240
$ws.Range("H86").Value = 224151.44
$ws.Range("I86").Value = 1910.6666
$ws.Range("J86").Value = 668633
$ws.Range("K86").Value = 1910.6666
$ws.Range("L86").Value = 668633
$ws.Range("M86").Value = -787.6666
$ws.Range("N86").Value = -670879
$ws.Range("H89").Value = 224151.44
$ws.Range("I89").Value = 1910.6666
$ws.Range("J89").Value = 668633
$ws.Range("K89").Value = 9553.333000000001
$ws.Range("L89").Value = 3343165
$ws.Range("M89").Value = -3937.333000000001
$ws.Range("N89").Value = -3354397
$ws.Range("H99").Value = 653.6667
$ws.Range("I99").Value = 675.375
$ws.Range("K99").Value = 675.375
$ws.Range("M99").Value = 822.625
$ws.Range("H134").Value = 6443.769
$ws.Range("I134").Value = 7148.8423
$ws.Range("K134").Value = 21446.5269
$ws.Range("M134").Value = -18911.5269
$ws.Range("H136").Value = 49800
$ws.Range("J136").Value = 49800
$ws.Range("L136").Value = 49800
$ws.Range("N136").Value = -60000

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 494.25
$ws.Range("I19").Value = 494.25
$ws.Range("K19").Value = 494.25
$ws.Range("M19").Value = -324.25
$ws.Range("H24").Value = 494.25
$ws.Range("I24").Value = 494.25
$ws.Range("K24").Value = 494.25
$ws.Range("M24").Value = -324.25
$ws.Range("H31").Value = 2747.6206
$ws.Range("I31").Value = 1413.5834
$ws.Range("J31").Value = 3689.2942
$ws.Range("K31").Value = 1413.5834
$ws.Range("L31").Value = 3689.2942
$ws.Range("M31").Value = -1118.5834
$ws.Range("N31").Value = -4279.2942
$ws.Range("H34").Value = 2747.6206
$ws.Range("I34").Value = 1413.5834
$ws.Range("J34").Value = 3689.2942
$ws.Range("K34").Value = 1413.5834
$ws.Range("L34").Value = 3689.2942
$ws.Range("M34").Value = -1211.5834
$ws.Range("N34").Value = -4093.2942
$ws.Range("H132").Value = 4349.6665
$ws.Range("I132").Value = 2356
$ws.Range("J132").Value = 5346.5
$ws.Range("K132").Value = 7068
$ws.Range("L132").Value = 16039.5
$ws.Range("M132").Value = -4538
$ws.Range("N132").Value = -21099.5
$ws.Range("H134").Value = 1387.1471
$ws.Range("I134").Value = 1288.68
$ws.Range("J134").Value = 1660.6666
$ws.Range("K134").Value = 3866.04
$ws.Range("L134").Value = 4981.9998
$ws.Range("M134").Value = -1331.04
$ws.Range("N134").Value = -10051.9998
$ws.Range("H141").Value = 64956.383
$ws.Range("J141").Value = 64369.418
$ws.Range("L141").Value = 64369.418
$ws.Range("N141").Value = -74729.41800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 318.16666
$ws.Range("I33").Value = 83.59999999999999
$ws.Range("J33").Value = 485.7143
$ws.Range("K33").Value = 501.6
$ws.Range("L33").Value = 2914.2858
$ws.Range("M33").Value = -218.6
$ws.Range("N33").Value = -3480.2858
$ws.Range("H80").Value = 2493.8
$ws.Range("J80").Value = 2493.8
$ws.Range("L80").Value = 7481.400000000001
$ws.Range("N80").Value = -9353.400000000001
$ws.Range("H83").Value = 2493.8
$ws.Range("J83").Value = 2493.8
$ws.Range("L83").Value = 22444.2
$ws.Range("N83").Value = -31804.2
$ws.Range("H118").Value = 1785.3334
$ws.Range("I118").Value = 1179
$ws.Range("J118").Value = 2998
$ws.Range("K118").Value = 3537
$ws.Range("L118").Value = 8994
$ws.Range("M118").Value = -2294
$ws.Range("N118").Value = -11480
$ws.Range("H125").Value = 5706.1816
$ws.Range("I125").Value = 1894.6666
$ws.Range("J125").Value = 10280
$ws.Range("K125").Value = 5683.9998
$ws.Range("L125").Value = 30840
$ws.Range("M125").Value = -763.9997999999996
$ws.Range("N125").Value = -40680
$ws.Range("H131").Value = 16093.705
$ws.Range("I131").Value = 451.33334
$ws.Range("J131").Value = 17800.145
$ws.Range("K131").Value = 1354.00002
$ws.Range("L131").Value = 53400.435
$ws.Range("M131").Value = 3685.99998
$ws.Range("N131").Value = -63480.435

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 117
$ws.Range("I2").Value = 48
$ws.Range("K2").Value = 48
$ws.Range("M2").Value = 65
$ws.Range("H46").Value = 22375
$ws.Range("J46").Value = 22375
$ws.Range("L46").Value = 22375
$ws.Range("N46").Value = -22687
$ws.Range("H80").Value = 2692.2354
$ws.Range("I80").Value = 2711.2
$ws.Range("J80").Value = 2550
$ws.Range("K80").Value = 2711.2
$ws.Range("L80").Value = 2550
$ws.Range("M80").Value = -1713.2
$ws.Range("N80").Value = -4546
$ws.Range("H83").Value = 2692.2354
$ws.Range("I83").Value = 2711.2
$ws.Range("J83").Value = 2550
$ws.Range("K83").Value = 13556
$ws.Range("L83").Value = 12750
$ws.Range("M83").Value = -8564
$ws.Range("N83").Value = -22734
$ws.Range("H113").Value = 1572.75
$ws.Range("I113").Value = 1220.75
$ws.Range("K113").Value = 1220.75
$ws.Range("M113").Value = 949.25
$ws.Range("H122").Value = 1518.8
$ws.Range("I122").Value = 1319
$ws.Range("J122").Value = 1818.5
$ws.Range("K122").Value = 3957
$ws.Range("L122").Value = 5455.5
$ws.Range("M122").Value = -1507
$ws.Range("N122").Value = -10355.5
$ws.Range("H126").Value = 2359515
$ws.Range("I126").Value = 2927124
$ws.Range("K126").Value = 8781372
$ws.Range("M126").Value = -8778902

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1279.6666
$ws.Range("I22").Value = 906.4286
$ws.Range("K22").Value = 906.4286
$ws.Range("M22").Value = -611.4286
$ws.Range("H27").Value = 1279.6666
$ws.Range("I27").Value = 906.4286
$ws.Range("K27").Value = 906.4286
$ws.Range("M27").Value = -799.4286
$ws.Range("H40").Value = 8674.963
$ws.Range("I40").Value = 7373.8184
$ws.Range("J40").Value = 14400
$ws.Range("K40").Value = 7373.8184
$ws.Range("L40").Value = 14400
$ws.Range("M40").Value = -7237.8184
$ws.Range("N40").Value = -14672

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 3200
$ws.Range("I23").Value = 3200
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 3200
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -2971
$ws.Range("N23").ClearContents()
$ws.Range("H81").Value = 874.875
$ws.Range("I81").Value = 874.875
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 1749.75
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -688.75
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 874.875
$ws.Range("I84").Value = 874.875
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 8748.75
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -3444.75
$ws.Range("N84").ClearContents()
$ws.Range("H132").Value = 1259.4348
$ws.Range("I132").Value = 1049.05
$ws.Range("J132").Value = 2662
$ws.Range("K132").Value = 3147.15
$ws.Range("L132").Value = 7986
$ws.Range("M132").Value = -617.1499999999996
$ws.Range("N132").Value = -13046
